$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - copy format from E1 (bold/bordered header style) then set value
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Timestamps for rows 2-51
$times = @(
    "2021-10-05 10:50:42.345708",
    "2021-10-05 10:50:42.345721",
    "2021-10-05 10:50:42.345725",
    "2021-10-05 10:50:42.345728",
    "2021-10-05 10:50:42.345732",
    "2021-10-05 10:50:42.345735",
    "2021-10-05 10:50:42.345738",
    "2021-10-05 10:50:42.345741",
    "2021-10-05 10:50:42.345744",
    "2021-10-05 10:50:42.345747",
    "2021-10-05 10:50:42.345750",
    "2021-10-05 10:50:42.345753",
    "2021-10-05 10:50:42.345756",
    "2021-10-05 10:50:42.345759",
    "2021-10-05 10:50:42.345762",
    "2021-10-05 10:50:42.345765",
    "2021-10-05 10:50:42.345768",
    "2021-10-05 10:50:42.345771",
    "2021-10-05 10:50:42.345774",
    "2021-10-05 10:50:42.345777",
    "2021-10-05 10:50:42.345780",
    "2021-10-05 10:50:42.345783",
    "2021-10-05 10:50:42.345786",
    "2021-10-05 10:50:42.345789",
    "2021-10-05 10:50:42.345792",
    "2021-10-05 10:50:42.345795",
    "2021-10-05 10:50:42.345798",
    "2021-10-05 10:50:42.345801",
    "2021-10-05 10:50:42.345804",
    "2021-10-05 10:50:42.345807",
    "2021-10-05 10:50:42.345810",
    "2021-10-05 10:50:42.345813",
    "2021-10-05 10:50:42.345816",
    "2021-10-05 10:50:42.345819",
    "2021-10-05 10:50:42.345822",
    "2021-10-05 10:50:42.345825",
    "2021-10-05 10:50:42.345828",
    "2021-10-05 10:50:42.345831",
    "2021-10-05 10:50:42.345834",
    "2021-10-05 10:50:42.345837",
    "2021-10-05 10:50:42.345840",
    "2021-10-05 10:50:42.345843",
    "2021-10-05 10:50:42.345847",
    "2021-10-05 10:50:42.345849",
    "2021-10-05 10:50:42.345852",
    "2021-10-05 10:50:42.345855",
    "2021-10-05 10:50:42.345858",
    "2021-10-05 10:50:42.345861",
    "2021-10-05 10:50:42.345864",
    "2021-10-05 10:50:42.345867"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
